# cond_pm1_train.xlsx: switch the recorded condition columns from the
# "angle difference" response coding to a location-judgment coding.
#
# New layout (row 1 header): SOA | stim1_ori | stim2_ori | stim1_c | stim2_c
# stim1_ori / stim2_ori hold "L" / "R" location labels instead of the old
# numeric stim1_c/stim2_c=1|0 columns, and the SOA-only condition table is
# collapsed down to the 4 distinct (SOA, stim1_ori, stim2_ori) rows that are
# actually used, each carrying the simplified -0.8 contrast value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row ---------------------------------------------------------
$ws.Range("A1").Value = "SOA"
$ws.Range("B1").Value = "stim1_ori"
$ws.Range("C1").Value = "stim2_ori"
$ws.Range("D1").Value = "stim1_c"
$ws.Range("E1").Value = "stim2_c"

# --- data rows -----------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "L"
$ws.Range("C2").Value = "L"
$ws.Range("D2").Value = -0.8
$ws.Range("E2").Value = -0.8

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = "L"
$ws.Range("C3").Value = "L"
$ws.Range("D3").Value = -0.8
$ws.Range("E3").Value = -0.8

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "L"
$ws.Range("C4").Value = "R"
$ws.Range("D4").Value = -0.8
$ws.Range("E4").Value = -0.8

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "L"
$ws.Range("C5").Value = "R"
$ws.Range("D5").Value = -0.8
$ws.Range("E5").Value = -0.8

# --- drop the old rows 6-13 (their data is gone in the new layout) ------
$ws.Range("A6:E13").ClearContents()

# The sheet keeps 5 trailing blank rows (6-10) at the original row height.
$ws.Rows.Item(6).RowHeight = 14.4
$ws.Rows.Item(7).RowHeight = 14.4
$ws.Rows.Item(8).RowHeight = 14.4
$ws.Rows.Item(9).RowHeight = 14.4
$ws.Rows.Item(10).RowHeight = 14.4
